$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44294
$ws.Cells.Item(2, 12).Value = "Especial"
$ws.Cells.Item(2, 13).Value = 200
$ws.Cells.Item(2, 14).Value = 14500
$ws.Cells.Item(2, 15).Value = 15000
$ws.Cells.Item(2, 16).Value = 14750
$ws.Cells.Item(2, 19).Value = 819

# Row 3
$ws.Cells.Item(3, 4).Value = 44294
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 240
$ws.Cells.Item(3, 14).Value = 12500
$ws.Cells.Item(3, 15).Value = 13000
$ws.Cells.Item(3, 16).Value = 12750
$ws.Cells.Item(3, 19).Value = 708

# Row 4
$ws.Cells.Item(4, 4).Value = 44294
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 240
$ws.Cells.Item(4, 14).Value = 10500
$ws.Cells.Item(4, 15).Value = 11000
$ws.Cells.Item(4, 16).Value = 10750
$ws.Cells.Item(4, 19).Value = 597

# Row 5
$ws.Cells.Item(5, 4).Value = 44606
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 240
$ws.Cells.Item(5, 14).Value = 11500
$ws.Cells.Item(5, 15).Value = 12000
$ws.Cells.Item(5, 16).Value = 11750
$ws.Cells.Item(5, 19).Value = 653

# Row 6
$ws.Cells.Item(6, 4).Value = 44606
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 240
$ws.Cells.Item(6, 14).Value = 9500
$ws.Cells.Item(6, 15).Value = 10000
$ws.Cells.Item(6, 16).Value = 9750
$ws.Cells.Item(6, 19).Value = 542

# Row 7
$ws.Cells.Item(7, 4).Value = 44631
$ws.Cells.Item(7, 12).Value = "Especial"
$ws.Cells.Item(7, 13).Value = 240
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 16000
$ws.Cells.Item(7, 16).Value = 15500
$ws.Cells.Item(7, 19).Value = 861

# Row 8
$ws.Cells.Item(8, 4).Value = 44631
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 248
$ws.Cells.Item(8, 14).Value = 12000
$ws.Cells.Item(8, 15).Value = 13000
$ws.Cells.Item(8, 16).Value = 12516
$ws.Cells.Item(8, 19).Value = 695

# Row 9
$ws.Cells.Item(9, 4).Value = 44631
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 9000
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 9500
$ws.Cells.Item(9, 19).Value = 528

# Row 10
$ws.Cells.Item(10, 4).Value = 44637
$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 200
$ws.Cells.Item(10, 14).Value = 14000
$ws.Cells.Item(10, 15).Value = 15000
$ws.Cells.Item(10, 16).Value = 14500
$ws.Cells.Item(10, 19).Value = 806

# Row 11
$ws.Cells.Item(11, 4).Value = 44637
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 240
$ws.Cells.Item(11, 14).Value = 10000
$ws.Cells.Item(11, 15).Value = 11000
$ws.Cells.Item(11, 16).Value = 10500
$ws.Cells.Item(11, 19).Value = 583

# Row 12
$ws.Cells.Item(12, 4).Value = 44636
$ws.Cells.Item(12, 12).Value = "Especial"
$ws.Cells.Item(12, 13).Value = 240
$ws.Cells.Item(12, 14).Value = 14000
$ws.Cells.Item(12, 15).Value = 15000
$ws.Cells.Item(12, 16).Value = 14500
$ws.Cells.Item(12, 19).Value = 806

# Row 13
$ws.Cells.Item(13, 4).Value = 44636
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 10000
$ws.Cells.Item(13, 15).Value = 11000
$ws.Cells.Item(13, 16).Value = 10500
$ws.Cells.Item(13, 19).Value = 583

# Row 14
$ws.Cells.Item(14, 4).Value = 44610
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 13000
$ws.Cells.Item(14, 15).Value = 14000
$ws.Cells.Item(14, 16).Value = 13500
$ws.Cells.Item(14, 19).Value = 750

# Row 15
$ws.Cells.Item(15, 4).Value = 44610
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 11000
$ws.Cells.Item(15, 15).Value = 12000
$ws.Cells.Item(15, 16).Value = 11500
$ws.Cells.Item(15, 19).Value = 639

# Row 16
$ws.Cells.Item(16, 4).Value = 44630
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 300
$ws.Cells.Item(16, 14).Value = 15000
$ws.Cells.Item(16, 15).Value = 16000
$ws.Cells.Item(16, 16).Value = 15500
$ws.Cells.Item(16, 19).Value = 861

# Row 17
$ws.Cells.Item(17, 4).Value = 44630
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 300
$ws.Cells.Item(17, 14).Value = 12000
$ws.Cells.Item(17, 15).Value = 13000
$ws.Cells.Item(17, 16).Value = 12500
$ws.Cells.Item(17, 19).Value = 694

# Row 18
$ws.Cells.Item(18, 4).Value = 44630
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 240
$ws.Cells.Item(18, 14).Value = 9000
$ws.Cells.Item(18, 15).Value = 10000
$ws.Cells.Item(18, 16).Value = 9500
$ws.Cells.Item(18, 19).Value = 528

# Row 19
$ws.Cells.Item(19, 4).Value = 44607
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 14).Value = 11000
$ws.Cells.Item(19, 15).Value = 12000
$ws.Cells.Item(19, 16).Value = 11500
$ws.Cells.Item(19, 19).Value = 639

# Row 20
$ws.Cells.Item(20, 4).Value = 44607
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 240
$ws.Cells.Item(20, 14).Value = 9000
$ws.Cells.Item(20, 15).Value = 10000
$ws.Cells.Item(20, 16).Value = 9500
$ws.Cells.Item(20, 19).Value = 528

# Row 21
$ws.Cells.Item(21, 4).Value = 44295
$ws.Cells.Item(21, 12).Value = "Especial"
$ws.Cells.Item(21, 13).Value = 200
$ws.Cells.Item(21, 14).Value = 14500
$ws.Cells.Item(21, 15).Value = 15000
$ws.Cells.Item(21, 16).Value = 14750
$ws.Cells.Item(21, 19).Value = 819

# Row 22
$ws.Cells.Item(22, 4).Value = 44295
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 12500
$ws.Cells.Item(22, 15).Value = 13000
$ws.Cells.Item(22, 16).Value = 12750
$ws.Cells.Item(22, 19).Value = 708

# Row 23
$ws.Cells.Item(23, 4).Value = 44295
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 240
$ws.Cells.Item(23, 14).Value = 10500
$ws.Cells.Item(23, 15).Value = 11000
$ws.Cells.Item(23, 16).Value = 10750
$ws.Cells.Item(23, 19).Value = 597

# Row 25
$ws.Cells.Item(25, 4).Value = 44634
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 200
$ws.Cells.Item(25, 14).Value = 14000
$ws.Cells.Item(25, 15).Value = 15000
$ws.Cells.Item(25, 16).Value = 14500
$ws.Cells.Item(25, 19).Value = 806

# Row 26
$ws.Cells.Item(26, 4).Value = 44634
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 10000
$ws.Cells.Item(26, 15).Value = 11000
$ws.Cells.Item(26, 16).Value = 10500
$ws.Cells.Item(26, 19).Value = 583

# Row 27
$ws.Cells.Item(27, 4).Value = 44603
$ws.Cells.Item(27, 12).Value = "Especial"
$ws.Cells.Item(27, 13).Value = 240
$ws.Cells.Item(27, 14).Value = 14500
$ws.Cells.Item(27, 15).Value = 15000
$ws.Cells.Item(27, 16).Value = 14750
$ws.Cells.Item(27, 19).Value = 819

# Row 28
$ws.Cells.Item(28, 4).Value = 44609
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 240
$ws.Cells.Item(28, 14).Value = 13000
$ws.Cells.Item(28, 15).Value = 14000
$ws.Cells.Item(28, 16).Value = 13500
$ws.Cells.Item(28, 19).Value = 750

# Row 29
$ws.Cells.Item(29, 4).Value = 44609
$ws.Cells.Item(29, 12).Value = "Segunda"
$ws.Cells.Item(29, 13).Value = 240
$ws.Cells.Item(29, 14).Value = 11000
$ws.Cells.Item(29, 15).Value = 12000
$ws.Cells.Item(29, 16).Value = 11500
$ws.Cells.Item(29, 19).Value = 639
